$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 21 data rows (rows 2-21) down by one row (to rows 3-22),
# working bottom-up so values are not clobbered before being read.
for ($r = 21; $r -ge 2; $r--) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r + 1, 1).Value2 = $a
    $ws.Cells.Item($r + 1, 2).Value2 = $b
    $ws.Cells.Item($r + 1, 3).Value2 = $c
}

# Write the new sample as the new first data row (row 2)
$ws.Cells.Item(2, 1).Value2 = -0.0467311926186084
$ws.Cells.Item(2, 2).Value2 = 0.0064140851609408
$ws.Cells.Item(2, 3).Value2 = -0.0200058370828628

# The shift above duplicated the former last data row into row 22; drop it so
# the sheet ends with header + 20 data rows (A1:C21).
$ws.Range("A22:C22").Clear()
